# Reduced number of items in study:
#  - Pre-test (col A) brought down to 10 items; Post-test (col C) has 25 items.
#  - None of the descriptions (col B) have corresponding images in the pre-test state
#    after this edit is fully reflected in the totals.
#  - Updates the A/B/C (Pre-test description / Pre-test image / Post-test image) grid
#    for rows 1-70, and the summary formulas in E2:G2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final A,B,C values for rows 1..70 (row 1 = header data row, matches existing layout)
$data = "0,1,0;1,0,1;0,0,0;1,0,1;0,0,0;1,0,1;0,1,1;1,0,0;0,0,0;1,0,0;0,0,1;0,0,0;0,0,0;0,1,1;0,0,0;0,0,0;1,0,1;0,0,0;0,0,1;1,0,0;1,0,1;0,0,0;0,0,1;0,0,0;0,0,1;1,0,0;0,0,0;0,0,0;1,0,1;0,1,1;0,0,0;0,0,0;0,0,1;0,0,0;0,0,0;0,0,1;0,0,0;0,0,0;0,0,1;0,0,0;0,0,0;0,0,1;0,0,0;0,0,0;0,0,1;0,1,1;0,0,1;0,0,0;0,0,0;0,0,1;0,0,0;0,0,0;0,0,0;0,0,1;0,0,0;0,0,0;0,0,1;0,0,0;0,0,0;0,1,0;0,0,0;0,0,0;0,1,0;0,0,0;0,1,0;0,0,1;0,1,1;0,0,0;0,0,0;0,1,0"

$rows = $data.Split(";")
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 1
    $vals = $rows[$i].Split(",")
    $ws.Cells.Item($r, 1).Value = [double]$vals[0]
    $ws.Cells.Item($r, 2).Value = [double]$vals[1]
    $ws.Cells.Item($r, 3).Value = [double]$vals[2]
}

# Summary formulas (post-test column B sum now also folds in the pre-test column A sum)
$ws.Range("E2").Formula = "=SUM(A1:A70)"
$ws.Range("F2").Formula = "=SUM(B1:B70) + SUM(A1:A70)"
$ws.Range("G2").Formula = "=SUM(C1:C70)"

# Column width tweak for column E (target stored width 12.15625 chars;
# ColumnWidth adds ~5/6 of a character of padding, so back that off here -
# the engine quantizes to 1/6-character steps, landing on the closest
# achievable stored width of ~12.1667)
$ws.Columns.Item(5).ColumnWidth = 11.322916666666666

# Restore the last active selection used when the file was saved
$ws.Range("C47").Select()

$wb.Save()
